$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Step 1: restructure rows via native worksheet row insert/delete ---
# Remove the old "Retrieve Policies based on Partner ID" row (old row 11)
$ws.Rows.Item(11).Delete()

# Insert 2 new rows before old "Map Policies to Partners" row (now at row 10)
# to make room for "Policy ID Generation" and "Policy ID Validation"
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# Insert 1 new row before "Distribution of Public Key to Partners" (now at row 16)
# to make room for "Partner Certificate Signing and RE-issueing"
$ws.Rows.Item(16).Insert()

# Resize the table to cover the new row count
$lo.Resize($ws.Range("B2:F20"))

# --- Step 2: write every cell value (B..F) for rows 1-20 to match target content ---
# Row 1
$ws.Range("B1").Value = "Partner Management (MISP and E-KYC/Auth Partners)"

# Row 2
$ws.Range("B2").Value = "Sr No."
$ws.Range("C2").Value = "Module"
$ws.Range("D2").Value = "Feature"
$ws.Range("E2").Value = "Acceptance Criteria"
$ws.Range("F2").Value = "Comments"

# Row 3
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Kernel"
$ws.Range("D3").Value = "MISP ID Generation"
$ws.Range("E3").Value = "1. Generate MISP ID as per below logic`na. MISP ID should be of 3 digits (Configurable)`nb. MISP ID should be generated sequentially`nc. MISP ID should be generated incrementally for every request"
$ws.Range("F3").Value = "Component already exist as TSP ID generator"

# Row 4
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Kernel"
$ws.Range("D4").Value = "MISP License Key Generation"
$ws.Range("E4").Value = "1. Generate a License Key as per below logic`na. License Key generation to follow random pattern`nb. License Key should be alphanumeric`nc. Length should be 8 digits (Configurable)`nd. Should be mapped to an expiry"

# Row 5
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "Kernel"
$ws.Range("D5").Value = "MISP License Key Pattern Validation"
$ws.Range("E5").Value = "1. Validate length of a License Key as configured and respond as mentioned below`na. If found valid, respond with `"Valid`"`nb. if found invalid, respond with `"Invalid`""

# Row 6
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "Admin"
$ws.Range("D6").Value = "MSIP License Key Expiry Validation"
$ws.Range("E6").Value = "1. Validate status of Lisence Key and respond as mentioned below`na. If found expired, respond with `"Your License Key is EXPIRED. Please regenrate a new License Key`"`nb. If found temporarily sespended, respond with `"Your License Key is temporarily SUSPENDED. Please contact MOSIP Administration`"`nc. If found permanently blocked, respond with `"Your License Key is BLOCKED. Please contact MOSIP Administration`""

# Row 7
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "Admin"
$ws.Range("D7").Value = "MISP Registration"
$ws.Range("E7").Value = "1. Receive request to register a MISP with follwing parameters`na. MISP Name`nb. MISP Contact Name`nc. MISP Phone`nd. MISP Email ID`n2. Issue and Map MISP ID`n3. Issue and Map Lisence Key`n4. Store the MISP in MOSIP"

# Row 8
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = "Kernel"
$ws.Range("D8").Value = "Partner ID Generation"
$ws.Range("E8").Value = "1. Generate Partner ID as per below logic`na. Partner ID should be of 4 digits (Configurable)`nb. Partner ID should be generated sequentially`nc. Partner ID should be generated incrementally for every request"

# Row 9
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = "Kernel"
$ws.Range("D9").Value = "Partner ID Validation"
$ws.Range("E9").Value = "1. Validate length of a Partner ID as configured and respond as mentioned below`na. If found valid, respond with `"Valid`"`nb. if found invalid, respond with `"Invalid`""

# Row 10
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "Kernel "
$ws.Range("D10").Value = "Policy ID Generation"
$ws.Range("E10").Value = "1. Generate Policy ID for following policies`na. OTP Trigger `nb. OTP Authentication`nc. Demo Authentication `nd. Biometric Authentication - FMR Data Match `ne. Biometric Authentication - IIR Data Match  `nf. Biometric Authentication - FID Data Match `ng. Static Pin Authentication`nh. eKYC - all combinations of eKYC demo fields `ni. Masked UIN`nj. UIN`n2. Generate Policy id as per below logic`na. Random ID generation`nb. Length should be 10 Digits (Configurable)"

# Row 11
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = "Kernel"
$ws.Range("D11").Value = "Policy ID Validation"
$ws.Range("E11").Value = "1. Validate length of a Policy ID as configured and respond as mentioned below`na. If found valid, respond with `"Valid`"`nb. if found invalid, respond with `"Invalid`""

# Row 12
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = "Admin"
$ws.Range("D12").Value = "Policy ID"
$ws.Range("E12").Value = "1. Receive request to retreive policies based on Partner ID and Policy ID`n2. Respond appropirately if Partner ID or Policy ID does not exist"

# Row 13
$ws.Range("B13").Value = 11
$ws.Range("C13").Value = "Admin"
$ws.Range("D13").Value = "Partner Registration"
$ws.Range("E13").Value = "1. Receive request to register a Partner with follwing parameters`na. Partner Name`nb. Partner Contact Name`nc. Partner Phone`nd. Partner Email ID`n2. Issue and Map Partner ID`n3. Map Policy ID to the Partner`na. Multiple Policies can be mapped to a Partner`nb. A Partner can have a policy for both Auth and E-KYC`n4. Store the Partner in MOSIP"

# Row 14
$ws.Range("B14").Value = 12
$ws.Range("C14").Value = "Admin"
$ws.Range("D14").Value = "MISP - Partner Mapping"
$ws.Range("E14").Value = "1. Receive a request to map MISP to a Partner with MISP ID and Partner ID as Input`n2. There can ve a many-to-mapping between MISPs and Partners"

# Row 15
$ws.Range("B15").Value = 13
$ws.Range("C15").Value = "Admin"
$ws.Range("D15").Value = "Partner Certiicate Validation"
$ws.Range("E15").Value = "1. Receive certificate from Partner`n2. Verify CA Authority of the certificate"

# Row 16
$ws.Range("B16").Value = 14
$ws.Range("C16").Value = "Admin"
$ws.Range("D16").Value = "Partner Certificate Signing and RE-issueing"
$ws.Range("E16").Value = "1. Receive certificate from Partner during Partner Registration`n2. Sign the Partner Certificate with MOSIP Private Key and issue a certificate chain`n3. Re-issue certficate back to the Partner`n4. Private key to change priodically as per the Key Rotation Policy set by admin"

# Row 17
$ws.Range("B17").Value = 15
$ws.Range("C17").Value = "Admin"
$ws.Range("D17").Value = "Distribution of Public Key to Partners"
$ws.Range("E17").Value = "1. Distribute Public Key to Partners correspinding to the Private Key used to signed the Certificate`n2. Public key needs to be distributed priodically whenever the Private Key is rotated"

# Row 18
$ws.Range("B18").Value = 16
$ws.Range("C18").Value = "Admin"
$ws.Range("D18").Value = "Device Registration"
$ws.Range("E18").Value = "TBD"
$ws.Range("F18").Value = "Yet to analyzed"

# Row 19
$ws.Range("B19").Value = 17
$ws.Range("C19").Value = "Admin"
$ws.Range("D19").Value = "Device Provider Registration"
$ws.Range("E19").Value = "TBD"
$ws.Range("F19").Value = "Yet to analyzed"

# Row 20
$ws.Range("B20").Value = 18
$ws.Range("C20").Value = "Admin"
$ws.Range("D20").Value = "RD Service Registration"
$ws.Range("E20").Value = "TBD"
$ws.Range("F20").Value = "Yet to analyzed"

# --- Step 3: row heights ---
$ws.Rows.Item(1).RowHeight = 21.5
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 58
$ws.Rows.Item(4).RowHeight = 72.5
$ws.Rows.Item(5).RowHeight = 58
$ws.Rows.Item(6).RowHeight = 101.5
$ws.Rows.Item(7).RowHeight = 116
$ws.Rows.Item(8).RowHeight = 58
$ws.Rows.Item(9).RowHeight = 58
$ws.Rows.Item(10).RowHeight = 203
$ws.Rows.Item(11).RowHeight = 58
$ws.Rows.Item(12).RowHeight = 43.5
$ws.Rows.Item(13).RowHeight = 145
$ws.Rows.Item(14).RowHeight = 43.5
$ws.Rows.Item(15).RowHeight = 29
$ws.Rows.Item(16).RowHeight = 87
$ws.Rows.Item(17).RowHeight = 58

# --- Step 4: special style for E16 (Partner Certificate Signing AC cell): left align + wrap ---
$e16 = $ws.Range("E16")
$e16.WrapText = $true
$e16.HorizontalAlignment = -4131  # xlLeft

# --- Step 5: selection ---
$ws.Range("E7").Select()

Write-Host "done"